$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Final")

# Update student count (A1): 6 -> 7
$ws.Range("A1").Value = 7

# Update K column quiz/lab scores for specific rows (0 -> new values)
$ws.Range("K3").Value = 2
$ws.Range("K4").Value = 2
$ws.Range("K5").Value = 2
$ws.Range("K9").Value = 1.7
$ws.Range("K17").Value = 1.9

# Update the active cell selection on the Final sheet (J18 -> K18)
$ws.Activate()
$ws.Range("K18").Select()
